$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.817.02'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '2.644.68'
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''537.07'
$ws.Range("D6").Value = '''144.00'
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("D9").Value = '''6.55'
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("E10").Value = '  +2.14%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").Value = '3.103.66'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").Value = '59.725.39'
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("D15").Value = '''20.92'
$ws.Range("E15").Value = '  +2.12%  '
$ws.Range("D16").Value = '2.645.00'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("D17").Value = '''0.0000134'
$ws.Range("E17").Value = '  +1.44%  '
$ws.Range("D18").Value = '''342.00'
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("D20").Value = '''10.21'
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").Value = '''6.39'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = '''1.00'
$ws.Range("D23").Value = '''67.60'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  +2.56%  '
$ws.Range("D28").Value = '0.0₃0752'
$ws.Range("E28").Value = '  +4.90%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '''1.66'
$ws.Range("E30").Value = '  +3.84%  '
$ws.Range("D32").Value = '''18.97'
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("D33").Value = '''150.88'
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("E35").Value = '  +2.05%  '
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").Value = '''0.835'
$ws.Range("E37").Value = '  +1.43%  '
$ws.Range("D38").Value = '''0.821'
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("D39").Value = '''288.74'
$ws.Range("E39").Value = '  +7.86%  '
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("D42").Value = '''0.602'
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = '''10.74'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").Value = '''0.0532'
$ws.Range("E45").Value = '  +3.56%  '
$ws.Range("D46").Value = '1.966.57'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("D48").Value = '''18.43'
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("D50").Value = '''111.02'
$ws.Range("E50").Value = '  -0.29%  '
